# Scheduled-runner refresh: pushes new currentAveragePrice* market data
# (columns H/I/J) into each profession sheet and rewrites the dependent
# Leve price/profit columns (K/L/M/N) that are kept as plain cached
# values (no formulas) in this workbook. A few rows also drop their HQ
# profit cell (N, or M+N) entirely where HQ crafting is no longer priced.

$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Clear-CellValue($ws, $row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 17 8 2051
Set-CellValue $ws 17 10 2338.75
Set-CellValue $ws 17 12 7016.25
Set-CellValue $ws 17 14 -7352.25

Set-CellValue $ws 19 8 604.58826
Set-CellValue $ws 19 10 548.1429000000001
Set-CellValue $ws 19 12 548.1429000000001
Set-CellValue $ws 19 14 -898.1429000000001

Set-CellValue $ws 32 8 7696892
Set-CellValue $ws 32 9 10600
Set-CellValue $ws 32 10 9094400
Set-CellValue $ws 32 11 10600
Set-CellValue $ws 32 12 9094400
Set-CellValue $ws 32 13 -10274
Set-CellValue $ws 32 14 -9095052

Set-CellValue $ws 113 8 4543.9
Set-CellValue $ws 113 10 5072.522
Set-CellValue $ws 113 12 5072.522
Set-CellValue $ws 113 14 -11580.522

Set-CellValue $ws 116 8 5016.3335
Set-CellValue $ws 116 9 4399.75
Set-CellValue $ws 116 11 4399.75
Set-CellValue $ws 116 13 -957.75

Set-CellValue $ws 125 8 1556.9166
Set-CellValue $ws 125 10 1632.4286
Set-CellValue $ws 125 12 14691.8574
Set-CellValue $ws 125 14 -19611.8574

Set-CellValue $ws 127 8 7807.44
Set-CellValue $ws 127 10 14962.909
Set-CellValue $ws 127 12 44888.727
Set-CellValue $ws 127 14 -54808.727

Set-CellValue $ws 135 8 20842978
Set-CellValue $ws 135 9 33338672
Set-CellValue $ws 135 11 300048048
Set-CellValue $ws 135 13 -300045513

$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 45 8 4459.7144
Set-CellValue $ws 45 9 2302.5557
Set-CellValue $ws 45 10 8342.6
Set-CellValue $ws 45 11 2302.5557
Set-CellValue $ws 45 12 8342.6
Set-CellValue $ws 45 13 -1925.5557
Set-CellValue $ws 45 14 -9096.6

Set-CellValue $ws 74 8 3093.8096
Set-CellValue $ws 74 9 3077.3684
Set-CellValue $ws 74 10 3250
Set-CellValue $ws 74 11 3077.3684
Set-CellValue $ws 74 12 3250
Set-CellValue $ws 74 13 -2203.3684
Set-CellValue $ws 74 14 -4998

Set-CellValue $ws 77 8 3093.8096
Set-CellValue $ws 77 9 3077.3684
Set-CellValue $ws 77 10 3250
Set-CellValue $ws 77 11 15386.842
Set-CellValue $ws 77 12 16250
Set-CellValue $ws 77 13 -11018.842
Set-CellValue $ws 77 14 -24986

Set-CellValue $ws 110 8 5137.2856
Set-CellValue $ws 110 9 5137.2856
Set-CellValue $ws 110 11 5137.2856
Set-CellValue $ws 110 13 -3092.2856

Set-CellValue $ws 122 8 3967.2666
Set-CellValue $ws 122 9 2286.4866
Set-CellValue $ws 122 11 6859.459800000001
Set-CellValue $ws 122 13 -4409.459800000001

Set-CellValue $ws 132 8 3376.3386
Set-CellValue $ws 132 9 3200.2942
Set-CellValue $ws 132 11 9600.882599999999
Set-CellValue $ws 132 13 -7070.882599999999

$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 20 8 36234.723
Set-CellValue $ws 20 10 126414.625
Set-CellValue $ws 20 12 126414.625
Set-CellValue $ws 20 14 -126908.625

Set-CellValue $ws 105 8 2997.4482
Set-CellValue $ws 105 10 4480.273
Set-CellValue $ws 105 12 4480.273
Set-CellValue $ws 105 14 -7974.273

Set-CellValue $ws 107 8 1701.7646
Set-CellValue $ws 107 9 1655.7097
Set-CellValue $ws 107 11 1655.7097
Set-CellValue $ws 107 13 264.2902999999999

$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 31 8 4310.778
Set-CellValue $ws 31 9 3457.6
Set-CellValue $ws 31 10 5377.25
Set-CellValue $ws 31 11 3457.6
Set-CellValue $ws 31 12 5377.25
Set-CellValue $ws 31 13 -3162.6
Set-CellValue $ws 31 14 -5967.25

Set-CellValue $ws 34 8 4310.778
Set-CellValue $ws 34 9 3457.6
Set-CellValue $ws 34 10 5377.25
Set-CellValue $ws 34 11 3457.6
Set-CellValue $ws 34 12 5377.25
Set-CellValue $ws 34 13 -3255.6
Set-CellValue $ws 34 14 -5781.25

Set-CellValue $ws 132 8 3283.7568
Set-CellValue $ws 132 9 2749.147
Set-CellValue $ws 132 11 8247.440999999999
Set-CellValue $ws 132 13 -5717.440999999999

Set-CellValue $ws 134 8 3590.8667
Set-CellValue $ws 134 9 2383.0908
Set-CellValue $ws 134 10 6912.25
Set-CellValue $ws 134 11 7149.2724
Set-CellValue $ws 134 12 20736.75
Set-CellValue $ws 134 13 -4614.2724
Set-CellValue $ws 134 14 -25806.75

$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws 3 8 3221.2856
Set-CellValue $ws 3 9 3221.2856
Set-CellValue $ws 3 10 0
Set-CellValue $ws 3 11 9663.856800000001
Set-CellValue $ws 3 12 0
Set-CellValue $ws 3 13 -9551.856800000001
Clear-CellValue $ws 3 14

Set-CellValue $ws 107 8 300.96155
Set-CellValue $ws 107 10 285.58334
Set-CellValue $ws 107 12 856.7500200000001
Set-CellValue $ws 107 14 -4696.75002

Set-CellValue $ws 111 8 2783
Set-CellValue $ws 111 9 449
Set-CellValue $ws 111 10 3950
Set-CellValue $ws 111 11 1347
Set-CellValue $ws 111 12 11850
Set-CellValue $ws 111 13 1720
Set-CellValue $ws 111 14 -17984

Set-CellValue $ws 116 8 1535.6666
Set-CellValue $ws 116 9 763.5
Set-CellValue $ws 116 11 2290.5
Set-CellValue $ws 116 13 1151.5

Set-CellValue $ws 120 8 0
Set-CellValue $ws 120 10 0
Set-CellValue $ws 120 12 0
Clear-CellValue $ws 120 14

Set-CellValue $ws 121 8 2231.889
Set-CellValue $ws 121 9 227.5
Set-CellValue $ws 121 10 2804.5715
Set-CellValue $ws 121 11 682.5
Set-CellValue $ws 121 12 8413.7145
Set-CellValue $ws 121 13 627.5
Set-CellValue $ws 121 14 -11033.7145

Set-CellValue $ws 131 8 3506.0356
Set-CellValue $ws 131 9 2987
Set-CellValue $ws 131 10 3794.389
Set-CellValue $ws 131 11 8961
Set-CellValue $ws 131 12 11383.167
Set-CellValue $ws 131 13 -3921
Set-CellValue $ws 131 14 -21463.167

Set-CellValue $ws 134 8 8502.846
Set-CellValue $ws 134 9 2256.1667
Set-CellValue $ws 134 11 6768.500100000001
Set-CellValue $ws 134 13 -1698.500100000001

Set-CellValue $ws 139 8 55559796
Set-CellValue $ws 139 9 62503650
Set-CellValue $ws 139 11 187510950
Set-CellValue $ws 139 13 -187505810

$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws 113 8 6345.9443
Set-CellValue $ws 113 9 3436.6667
Set-CellValue $ws 113 10 12164.5
Set-CellValue $ws 113 11 3436.6667
Set-CellValue $ws 113 12 12164.5
Set-CellValue $ws 113 13 -1266.6667
Set-CellValue $ws 113 14 -16504.5

Set-CellValue $ws 119 8 50204.4
Set-CellValue $ws 119 10 50204.4
Set-CellValue $ws 119 12 50204.4
Set-CellValue $ws 119 14 -59880.4

Set-CellValue $ws 126 8 5728.7144
Set-CellValue $ws 126 9 3867.3333
Set-CellValue $ws 126 11 11601.9999
Set-CellValue $ws 126 13 -9131.999899999999

$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 16 8 113.1875
Set-CellValue $ws 16 9 113.1875
Set-CellValue $ws 16 10 0
Set-CellValue $ws 16 11 113.1875
Set-CellValue $ws 16 12 0
Set-CellValue $ws 16 13 56.8125
Clear-CellValue $ws 16 14

Set-CellValue $ws 61 8 3250.138
Set-CellValue $ws 61 9 1655
Set-CellValue $ws 61 11 1655
Set-CellValue $ws 61 13 -1453

Set-CellValue $ws 113 8 3250.138
Set-CellValue $ws 113 9 1655
Set-CellValue $ws 113 11 1655
Set-CellValue $ws 113 13 515

Set-CellValue $ws 122 8 4356.9414
Set-CellValue $ws 122 9 3561.9644
Set-CellValue $ws 122 10 8066.8335
Set-CellValue $ws 122 11 10685.8932
Set-CellValue $ws 122 12 24200.5005
Set-CellValue $ws 122 13 -8235.893199999999
Set-CellValue $ws 122 14 -29100.5005

$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws 96 8 2293.111
Set-CellValue $ws 96 9 2152.4666
Set-CellValue $ws 96 10 2996.3333
Set-CellValue $ws 96 11 2152.4666
Set-CellValue $ws 96 12 2996.3333
Set-CellValue $ws 96 13 -779.4666000000002
Set-CellValue $ws 96 14 -5742.3333

Set-CellValue $ws 122 8 5021.069
Set-CellValue $ws 122 9 2580.6956
Set-CellValue $ws 122 11 7742.0868
Set-CellValue $ws 122 13 -5292.0868

Set-CellValue $ws 126 8 2256.5833
Set-CellValue $ws 126 9 2387.4
Set-CellValue $ws 126 10 1602.5
Set-CellValue $ws 126 11 7162.200000000001
Set-CellValue $ws 126 12 4807.5
Set-CellValue $ws 126 13 -4692.200000000001
Set-CellValue $ws 126 14 -9747.5

Set-CellValue $ws 132 8 1681.1177
Set-CellValue $ws 132 9 1286.1875
Set-CellValue $ws 132 11 3858.5625
Set-CellValue $ws 132 13 -1328.5625

Write-Output "Applied scheduled market-data refresh to all 8 sheets."
